# Apply the "through September 10" data refresh to the carjacking-by-
# neighborhood-by-month workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Sheet title / header text: "September 09" -> "September 10" ---
$wb.Worksheets.Item(1).Name = "Through 2021-09-10"
$ws.Cells.Item(1, 2).Value = "September 2021 (through September 10)"

# --- Cell value updates (row, col) -> new value ---
$updates = @(
    @(3,  2,  2),
    @(4,  2,  2),
    @(4,  11, 3),
    @(4,  38, 1),
    @(5,  2,  4),
    @(5,  29, 2),
    @(7,  11, 1),
    @(13, 20, 2),
    @(16, 38, 2),
    @(18, 38, 2),
    @(19, 2,  1),
    @(26, 2,  1),
    @(31, 11, 2),
    @(38, 20, 1),
    @(41, 47, 1),
    @(47, 56, 1),
    @(55, 2,  2),
    @(63, 47, 1),
    @(67, 11, 1),
    @(93, 2,  1)
)

foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}
